$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 1897667.4
$ws.Range("I38").Value = 2016235.2
$ws.Range("J38").Value = 580
$ws.Range("K38").Value = 6048705.6
$ws.Range("L38").Value = 1740
$ws.Range("M38").Value = -6048333.6
$ws.Range("N38").Value = -2484

# Row 62
$ws.Range("H62").Value = 1953.2222
$ws.Range("I62").Value = 1930
$ws.Range("J62").Value = 1999.6666
$ws.Range("K62").Value = 1930
$ws.Range("L62").Value = 1999.6666
$ws.Range("M62").Value = -1306
$ws.Range("N62").Value = -3247.6666

# Row 65
$ws.Range("H65").Value = 1953.2222
$ws.Range("I65").Value = 1930
$ws.Range("J65").Value = 1999.6666
$ws.Range("K65").Value = 9650
$ws.Range("L65").Value = 9998.333000000001
$ws.Range("M65").Value = -6530
$ws.Range("N65").Value = -16238.333

# Row 99
$ws.Range("H99").Value = 368
$ws.Range("I99").Value = 228.8
$ws.Range("J99").Value = 600
$ws.Range("K99").Value = 686.4000000000001
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = 811.5999999999999
$ws.Range("N99").Value = -4796

# Row 101
$ws.Range("H101").Value = 454
$ws.Range("I101").Value = 454
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 1362
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 260
$ws.Range("N101").ClearContents()

# Row 115
$ws.Range("H115").Value = 887
$ws.Range("I115").Value = 792.5
$ws.Range("J115").Value = 950
$ws.Range("K115").Value = 2377.5
$ws.Range("L115").Value = 2850
$ws.Range("M115").Value = -810.5

# Row 118
$ws.Range("H118").Value = 1098
$ws.Range("I118").Value = 490
$ws.Range("J118").Value = 1250
$ws.Range("K118").Value = 1470
$ws.Range("L118").Value = 3750
$ws.Range("M118").Value = 187
$ws.Range("N118").Value = -7064

# Row 127
$ws.Range("H127").Value = 1714.4324
$ws.Range("I127").Value = 524.3
$ws.Range("J127").Value = 2155.2222
$ws.Range("K127").Value = 1572.9
$ws.Range("L127").Value = 6465.6666
$ws.Range("M127").Value = 3387.1
$ws.Range("N127").Value = -16385.6666

# Row 129
$ws.Range("H129").Value = 855.5
$ws.Range("I129").Value = 574.6875
$ws.Range("J129").Value = 980.30554
$ws.Range("K129").Value = 1724.0625
$ws.Range("L129").Value = 2940.91662
$ws.Range("M129").Value = 3275.9375
$ws.Range("N129").Value = -12940.91662

# Row 138
$ws.Range("H138").Value = 3048
$ws.Range("I138").Value = 2539
$ws.Range("J138").Value = 3159.0544
$ws.Range("K138").Value = 7617
$ws.Range("L138").Value = 9477.163199999999
$ws.Range("M138").Value = -2477
$ws.Range("N138").Value = -19757.1632

$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 7166.5
$ws.Range("I3").Value = 2999
$ws.Range("J3").Value = 8000
$ws.Range("K3").Value = 2999
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = -2884
$ws.Range("N3").Value = -8230

# Row 11
$ws.Range("H11").Value = 4566.6665
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 4566.6665
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 4566.6665
$ws.Range("N11").Value = -4854.6665

# Row 118
$ws.Range("H118").Value = 35950
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 35950
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 35950
$ws.Range("N118").Value = -39264

$ws = $wb.Worksheets.Item("BSM")
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 19
$ws.Range("H19").Value = 24900
$ws.Range("I19").Value = 1800
$ws.Range("J19").Value = 48000
$ws.Range("K19").Value = 1800
$ws.Range("L19").Value = 48000
$ws.Range("M19").Value = -1630
$ws.Range("N19").Value = -48340

# Row 24
$ws.Range("H24").Value = 24900
$ws.Range("I24").Value = 1800
$ws.Range("J24").Value = 48000
$ws.Range("K24").Value = 1800
$ws.Range("L24").Value = 48000
$ws.Range("M24").Value = -1630
$ws.Range("N24").Value = -48340

# Row 98
$ws.Range("H98").Value = 48000
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 48000
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 48000
$ws.Range("N98").Value = -52492

$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 500.75
$ws.Range("I46").Value = 500.75
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1502.25
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1411.25
$ws.Range("N46").ClearContents()

# Row 51
$ws.Range("H51").Value = 2147.5
$ws.Range("I51").Value = 2147.5
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 6442.5
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -5982.5
$ws.Range("N51").ClearContents()

# Row 58
$ws.Range("H58").Value = 850
$ws.Range("I58").Value = 850
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2550
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2422

# Row 98
$ws.Range("H98").Value = 69539.75
$ws.Range("I98").Value = 525.75
$ws.Range("J98").Value = 92544.414
$ws.Range("K98").Value = 1577.25
$ws.Range("L98").Value = 277633.242
$ws.Range("M98").Value = -79.25
$ws.Range("N98").Value = -280629.242

# Row 114
$ws.Range("H114").Value = 849.7143
$ws.Range("I114").Value = 66.59999999999999
$ws.Range("J114").Value = 1284.7778
$ws.Range("K114").Value = 199.8
$ws.Range("L114").Value = 3854.3334
$ws.Range("M114").Value = 3054.2
$ws.Range("N114").Value = -10362.3334

# Row 119
$ws.Range("H119").Value = 4600.385
$ws.Range("I119").Value = 621
$ws.Range("J119").Value = 7087.5
$ws.Range("K119").Value = 1863
$ws.Range("L119").Value = 21262.5
$ws.Range("M119").Value = 2975
$ws.Range("N119").Value = -30938.5

# Row 120
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("N120").ClearContents()

# Row 129
$ws.Range("H129").Value = 11906226
$ws.Range("I129").Value = 50000620
$ws.Range("J129").Value = 1727.5625
$ws.Range("K129").Value = 150001860
$ws.Range("L129").Value = 5182.6875
$ws.Range("M129").Value = -149996860
$ws.Range("N129").Value = -15182.6875

# Row 131
$ws.Range("H131").Value = 866.8523
$ws.Range("I131").Value = 432.85715
$ws.Range("J131").Value = 904.35803
$ws.Range("K131").Value = 1298.57145
$ws.Range("L131").Value = 2713.07409
$ws.Range("M131").Value = 3741.42855
$ws.Range("N131").Value = -12793.07409

# Row 137
$ws.Range("H137").Value = 18521140
$ws.Range("I137").Value = 2023.75
$ws.Range("J137").Value = 33336434
$ws.Range("K137").Value = 6071.25
$ws.Range("L137").Value = 100009302
$ws.Range("M137").Value = -971.25
$ws.Range("N137").Value = -100019502

# Row 140
$ws.Range("H140").Value = 1271.579
$ws.Range("I140").Value = 1064.4445
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 3193.3335
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = 1986.6665

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 2933.3333
$ws.Range("I3").Value = 1300
$ws.Range("J3").Value = 3750
$ws.Range("K3").Value = 1300
$ws.Range("L3").Value = 3750
$ws.Range("M3").Value = -1188
$ws.Range("N3").Value = -3974

# Row 15
$ws.Range("H15").Value = 2933.3333
$ws.Range("I15").Value = 1300
$ws.Range("J15").Value = 3750
$ws.Range("K15").Value = 1300
$ws.Range("L15").Value = 3750
$ws.Range("M15").Value = -1130
$ws.Range("N15").Value = -4090

# Row 55
$ws.Range("H55").Value = 534.94446
$ws.Range("I55").Value = 263.875
$ws.Range("J55").Value = 751.8
$ws.Range("K55").Value = 263.875
$ws.Range("L55").Value = 751.8
$ws.Range("M55").Value = -90.875
$ws.Range("N55").Value = -1097.8

# Row 93
$ws.Range("H93").Value = 1672.826
$ws.Range("I93").Value = 1343.1666
$ws.Range("J93").Value = 2859.6
$ws.Range("K93").Value = 1343.1666
$ws.Range("L93").Value = 2859.6
$ws.Range("M93").Value = -95.16660000000002
$ws.Range("N93").Value = -5355.6

$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 37980
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 37980
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 37980
$ws.Range("N16").Value = -38564
